# Update existing data rows (2-6) with new computed values, then append two
# new iteration rows (7-8) to extend the convergence table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows 2-6: overwrite xi / f(xi) / Error columns -------------------
# Force text storage (these are numeric-looking strings stored as text in
# the source data, same as the rest of the sheet) without leaving a
# lingering custom cell style behind.
$dataRange = $ws.Range("B2:D6")
$dataRange.NumberFormat = "@"

$ws.Range("B2").Value = "1.45098039215686"
$ws.Range("C2").Value = "47.0"
$ws.Range("D2").Value = "6.45098039215686"

$ws.Range("B3").Value = "1.04577347143334"
$ws.Range("C3").Value = "-1.69857747020377"
$ws.Range("D3").Value = "0.405206920723518"

$ws.Range("B4").Value = "0.688595493685113"
$ws.Range("C4").Value = "-1.08945173217963"
$ws.Range("D4").Value = "0.357177977748231"

$ws.Range("B5").Value = "0.589804903444715"
$ws.Range("C5").Value = "-0.280218220817008"
$ws.Range("D5").Value = "0.0987905902403977"

$ws.Range("B6").Value = "0.585792163076594"
$ws.Range("C6").Value = "-0.0113497896514318"
$ws.Range("D6").Value = "0.0040127403681217"

$dataRange.Style = "Normal"

# --- Rows 7-8: new iterations (table grows from A1:D6 to A1:D8) -------
$newRows = $ws.Range("A7:D8")
$newRows.NumberFormat = "@"

$ws.Range("A7").Value = "6"
$ws.Range("B7").Value = "0.585786437638495"
$ws.Range("C7").Value = "-1.61939844204328e-05"
$ws.Range("D7").Value = "5.72543809895798e-06"

$ws.Range("A8").Value = "7"
$ws.Range("B8").Value = "0.585786437626905"
$ws.Range("C8").Value = "-3.2781058975394394e-11"
$ws.Range("D8").Value = "1.1589840198666899e-11"

$newRows.Style = "Normal"
